# #5: property aircraft done
# The 建物 (Building) sheet rows were mistakenly tagged with the "land"
# property_category, and the 汽車 (Car) sheet rows were mistakenly tagged
# with the "land" category as well. Fix both to their correct values.

$wb = $excel.ActiveWorkbook

# 建物 (Building) sheet: property_category column (I) for the two data
# rows should read "building" instead of "land".
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"

# 汽車 (Car) sheet: category column (H) for the three data rows should
# read "car" instead of "land".
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
$wsCar.Range("H3").Value = "car"
$wsCar.Range("H4").Value = "car"
